$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear segment/time labels, update scale
$ws.Range("A2:B2").ClearContents()
$ws.Range("C2").Value = 1

# Row 3: rename "pi pulse" -> "π pulse", update scale
$ws.Range("A3").Value = "π pulse"
$ws.Range("C3").Value = 1

# Row 4: update scale
$ws.Range("C4").Value = 4

# Row 5: update scale
$ws.Range("C5").Value = 3

# Row 6: clear segment/time labels, update scale
$ws.Range("A6:B6").ClearContents()
$ws.Range("C6").Value = 1

# Row 7: update scale
$ws.Range("C7").Value = 0.5

# Row 8: update scale
$ws.Range("C8").Value = 4

# Row 9: clear segment/time labels and old param2 columns, update scale, add new param1 columns
$ws.Range("A9:B9").ClearContents()
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = "linear"
$ws.Range("F9:G9").ClearContents()

# Row 10: rename "SSI" -> "sweep", update time, scale, and add param2 columns
$ws.Range("A10").Value = "sweep"
$ws.Range("B10").Value = "9ms"
$ws.Range("C10").Value = 2.5
$ws.Range("F10").Value = 209
$ws.Range("G10").Value = "linear"

# Row 11 (new)
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = "linear"

# Row 12 (new)
$ws.Range("C12").Value = 1

# Selection moved to B10 as in the saved file
$ws.Range("B10").Select()
